$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D, L, M, N, O, P, S per row, after the
# weekly re-consolidation reshuffled which daily record sits in which row.
$rows = @{
    2 = @{ D = 44630; L = "Especial"; M = 300; N = 15000; O = 16000; P = 15500; S = 861 }
    3 = @{ D = 44630; L = "Primera"; M = 300; N = 12000; O = 13000; P = 12500; S = 694 }
    4 = @{ D = 44630; L = "Segunda"; M = 240; N = 9000; O = 10000; P = 9500; S = 528 }
    5 = @{ D = 44637; L = "Especial"; M = 200; N = 14000; O = 15000; P = 14500; S = 806 }
    6 = @{ D = 44637; L = "Primera"; M = 240; N = 10000; O = 11000; P = 10500; S = 583 }
    7 = @{ D = 44609; L = "Primera"; M = 240; N = 13000; O = 14000; P = 13500; S = 750 }
    8 = @{ D = 44609; L = "Segunda"; M = 240; N = 11000; O = 12000; P = 11500; S = 639 }
    9 = @{ D = 44295; L = "Especial"; M = 200; N = 14500; O = 15000; P = 14750; S = 819 }
    10 = @{ D = 44295; L = "Primera"; M = 200; N = 12500; O = 13000; P = 12750; S = 708 }
    11 = @{ D = 44295; L = "Segunda"; M = 240; N = 10500; O = 11000; P = 10750; S = 597 }
    12 = @{ D = 44606; L = "Primera"; M = 240; N = 11500; O = 12000; P = 11750; S = 653 }
    13 = @{ D = 44606; L = "Segunda"; M = 240; N = 9500; O = 10000; P = 9750; S = 542 }
    14 = @{ D = 44636; L = "Especial"; M = 240; N = 14000; O = 15000; P = 14500; S = 806 }
    15 = @{ D = 44636; L = "Primera"; M = 200; N = 10000; O = 11000; P = 10500; S = 583 }
    16 = @{ D = 44610; L = "Primera"; M = 200; N = 13000; O = 14000; P = 13500; S = 750 }
    17 = @{ D = 44610; L = "Segunda"; M = 200; N = 11000; O = 12000; P = 11500; S = 639 }
    18 = @{ D = 44603; L = "Especial"; M = 240; N = 14500; O = 15000; P = 14750; S = 819 }
    19 = @{ D = 44631; L = "Especial"; M = 240; N = 15000; O = 16000; P = 15500; S = 861 }
    20 = @{ D = 44631; L = "Primera"; M = 248; N = 12000; O = 13000; P = 12516; S = 695 }
    21 = @{ D = 44631; L = "Segunda"; M = 200; N = 9000; O = 10000; P = 9500; S = 528 }
    22 = @{ D = 44634; L = "Especial"; M = 200; N = 14000; O = 15000; P = 14500; S = 806 }
    23 = @{ D = 44634; L = "Primera"; M = 200; N = 10000; O = 11000; P = 10500; S = 583 }
    24 = @{ D = 44595; L = "Primera"; M = 200; N = 15500; O = 16000; P = 15750; S = 875 }
    25 = @{ D = 44294; L = "Especial"; M = 200; N = 14500; O = 15000; P = 14750; S = 819 }
    26 = @{ D = 44294; L = "Primera"; M = 240; N = 12500; O = 13000; P = 12750; S = 708 }
    27 = @{ D = 44294; L = "Segunda"; M = 240; N = 10500; O = 11000; P = 10750; S = 597 }
    28 = @{ D = 44607; L = "Primera"; M = 300; N = 11000; O = 12000; P = 11500; S = 639 }
    29 = @{ D = 44607; L = "Segunda"; M = 240; N = 9000; O = 10000; P = 9500; S = 528 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $v.D
    $ws.Cells.Item($r, 12).Value = $v.L
    $ws.Cells.Item($r, 13).Value = $v.M
    $ws.Cells.Item($r, 14).Value = $v.N
    $ws.Cells.Item($r, 15).Value = $v.O
    $ws.Cells.Item($r, 16).Value = $v.P
    $ws.Cells.Item($r, 19).Value = $v.S
}
